# Apply updated dSF (column F) values per repulled data / mean calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -3
    4  = -7
    6  = -1
    7  = 2
    8  = 5
    9  = 11
    10 = 1
    11 = -3
    12 = -5
    13 = -1
    14 = -2
    15 = 1
    19 = -3
    20 = -4
    21 = 3
    22 = -5
    23 = -3
    24 = 1
    25 = -1
    26 = -3
    27 = -3
    28 = -4
    29 = -2
    30 = -3
    31 = -2
    32 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
